# Appends newly-screened patient rows (840-858) to Sheet1, mirroring the
# "exports all new patients now" export. Columns: A name, B datescreen,
# C race_txt, D gender_txt, E age_at_encounter, F dob, G insurance, H zipcode.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(840, "Anaudia Johnson",      45221, "Black, Not Hispanic", "Woman", 24.28523515198806,  36351,   "NA", $null),
    @(841, "Barney rogers",        45223, "Black, Not Hispanic", "Man",   999,                -306287,  "NA", "27503"),
    @(842, "Denisha Smith",        45224, "Black, Not Hispanic", "Woman", 24.95864619624861,  36108,   "NA", "27707"),
    @(843, "Carla Green",          45227, "Black, Not Hispanic", "Woman", 34.36073293770577,  32677,   "NA", "27707"),
    @(844, "Emily Harrington",     45230, "White, Not Hispanic", "Woman", 22.7683205906578,   36914,   "NA", "27519"),
    @(845, "Deborah Chapman",      45230, "White, Not Hispanic", "Woman", 56.75396026840615,  24501,   "NA", "27701"),
    @(846, "Cynthia Parrish Fox",  45230, "Other",                "Woman", 67.29501632477053,  20651,   "NA", "27707"),
    @(847, "Ja'Naise Allison",     45232, "Black, Not Hispanic", "Woman", 0.7226933703863415, 44968,   "NA", "27701"),
    @(848, "Sandra McGhee-Bureh",  45235, "Black, Not Hispanic", "Woman", 67.61797527213655,  20538,   "NA", "27603"),
    @(849, "Jean McGhee",          45235, $null,                  "Woman", 0,                   45235,   "NA", "27573"),
    @(850, "Deborah E. Baker",     45235, "White, Not Hispanic", "Woman", 67.20466539353991,  20689,   "NA", "27217"),
    @(851, "Juan de Dios Argueta", 45236, "Other",                "Man",   24.33189365056549,  36349,   "NA", "27701"),
    @(852, "Waleed Razzaq",        45238, "White, Not Hispanic", "Man",   60.74320485704703,  23052,   "NA", "27893"),
    @(853, "Sandra McFadgir",      45238, "Black, Not Hispanic", "Woman", 61.95609766114294,  22609,   "NA", "27703"),
    @(854, "Vickie Bailey",        45238, "White, Not Hispanic", "Woman", 58.93618623243461,  23712,   "NA", $null),
    @(855, "Karen Richard",        45242, "White, Not Hispanic", "Woman", 67.82617028412631,  20469,   "NA", "27574"),
    @(856, "Brian Long",           45244, "Black, Not Hispanic", "Man",   53.87927199052684,  25565,   "NA", "27703"),
    @(857, "linda Hunter",         45244, "Black, Not Hispanic", "Woman", 35.96788435080803,  32107,   "NA", "27513"),
    @(858, "Bryce Bates",          45244, "White, Not Hispanic", "Man",   $null,               $null,   "NA", "27705")
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]

    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }

    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }

    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }

    if ($row[5] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[5] }

    if ($row[6] -ne $null) { $ws.Cells.Item($r, 6).Value = $row[6] }

    if ($row[7] -ne $null) { $ws.Cells.Item($r, 7).Value = $row[7] }

    if ($row[8] -ne $null) {
        $cell = $ws.Cells.Item($r, 8)
        $cell.NumberFormat = "@"
        $cell.Value = $row[8]
    }
}
